$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1) "Eccezioni" table cell: remove the explanatory sentence, leaving the
#    paragraph empty (its pPr / bold formatting definition stays intact).
# -----------------------------------------------------------------------
$oldEccezioni = "Non esistono aziende aderenti al sistema nella città in cui il Cliente o il Guest hanno cercato."
$foundEcc = $d.Content.Find.Execute($oldEccezioni, $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
Write-Host "Eccezioni text removed: $foundEcc"

# -----------------------------------------------------------------------
# 2) "Scenario" paragraph: rewrite the description of how the System
#    behaves, then append a brand new sentence after the bookmark.
# -----------------------------------------------------------------------
$oldScenario = "gestisce in caso ci siano o meno aziende nella città inserita, il Sistema fa visualizzare l’elenco di aziende nella città inserita."
$newScenario = "mostra le aziende nella città inserita e in caso non ci siano aziende in quella zona "
$foundScenario = $d.Content.Find.Execute($oldScenario, $true, $false, $false, $false, $false, $true, 1, $false, $newScenario, 2)
Write-Host "Scenario text rewritten: $foundScenario"

# Append the new closing sentence "non mostra nulla." right after the
# existing paragraph content (i.e. after the _GoBack bookmark), bold +
# bold-complex-script, matching the rest of the paragraph's formatting.
$scenarioParaIndex = $d.Paragraphs.Count
$scenarioPara = $d.Paragraphs.Item($scenarioParaIndex)
$scenarioPara.Range.InsertAfter("non mostra nulla.")

$newRunRange = $d.Content
$foundNewRun = $newRunRange.Find.Execute("non mostra nulla.")
Write-Host "New closing sentence inserted: $foundNewRun"
$newRunRange.Font.Bold = 1
$newRunRange.Font.BoldBi = 1
